$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fold_1")
$ws.Range("B2").Value = 22.99835
$ws.Range("C2").Value = 9.009249999999998
$ws.Range("D2").Value = 197.02425
$ws.Range("F2").Value = 6495.752149999998
$ws.Range("G2").Value = 5362.409900000001
$ws.Range("H2").Value = 1133.3425
$ws.Range("I2").Value = 1133.3425
$ws.Range("K2").Value = 5983.008599999999
$ws.Range("L2").Value = 5362.422
$ws.Range("M2").Value = 620.5866
$ws.Range("N2").Value = 620.5866
$ws.Range("B3").Value = 22.867
$ws.Range("C3").Value = 9.141
$ws.Range("D3").Value = 196.595
$ws.Range("F3").Value = 6558.194
$ws.Range("G3").Value = 5349.126
$ws.Range("H3").Value = 1209.068
$ws.Range("I3").Value = 1209.068
$ws.Range("K3").Value = 5979.8758
$ws.Range("L3").Value = 5349.074
$ws.Range("M3").Value = 630.8018
$ws.Range("N3").Value = 630.8018
$ws.Range("B4").Value = 26.357
$ws.Range("C4").Value = 10.383
$ws.Range("D4").Value = 209.797
$ws.Range("F4").Value = 7515.411
$ws.Range("G4").Value = 5900.331
$ws.Range("H4").Value = 1615.08
$ws.Range("I4").Value = 1615.08
$ws.Range("K4").Value = 6187.2038
$ws.Range("L4").Value = 5900.331
$ws.Range("M4").Value = 286.8728
$ws.Range("N4").Value = 286.8728

$ws = $wb.Worksheets.Item("Fold_2")
$ws.Range("B2").Value = 22.89015
$ws.Range("C2").Value = 9.385699999999998
$ws.Range("D2").Value = 194.3002
$ws.Range("F2").Value = 6360.090999999999
$ws.Range("G2").Value = 5343.9704
$ws.Range("H2").Value = 1016.1209
$ws.Range("I2").Value = 1016.1209
$ws.Range("K2").Value = 6495.248399999999
$ws.Range("L2").Value = 5343.98
$ws.Range("M2").Value = 1151.2684
$ws.Range("N2").Value = 1151.2684
$ws.Range("B3").Value = 22.814
$ws.Range("C3").Value = 9.388999999999999
$ws.Range("D3").Value = 194.162
$ws.Range("F3").Value = 6428.854
$ws.Range("G3").Value = 5334.268
$ws.Range("H3").Value = 1094.586
$ws.Range("I3").Value = 1094.586
$ws.Range("K3").Value = 6494.179
$ws.Range("L3").Value = 5334.25
$ws.Range("M3").Value = 1159.929
$ws.Range("N3").Value = 1159.929
$ws.Range("B4").Value = 26.357
$ws.Range("C4").Value = 10.383
$ws.Range("D4").Value = 209.797
$ws.Range("F4").Value = 7515.411
$ws.Range("G4").Value = 5900.331
$ws.Range("H4").Value = 1615.08
$ws.Range("I4").Value = 1615.08
$ws.Range("K4").Value = 6634.525200000001
$ws.Range("L4").Value = 5900.331
$ws.Range("M4").Value = 734.1942
$ws.Range("N4").Value = 734.1942

$ws = $wb.Worksheets.Item("Fold_3")
$ws.Range("B2").Value = 22.58955
$ws.Range("C2").Value = 9.17435
$ws.Range("D2").Value = 191.23895
$ws.Range("F2").Value = 6330.07
$ws.Range("G2").Value = 5278.01935
$ws.Range("H2").Value = 1052.0509
$ws.Range("I2").Value = 1052.0509
$ws.Range("K2").Value = 6616.836
$ws.Range("L2").Value = 5278.029
$ws.Range("M2").Value = 1338.807
$ws.Range("N2").Value = 1338.807
$ws.Range("B3").Value = 22.705
$ws.Range("C3").Value = 9.209
$ws.Range("D3").Value = 188.364
$ws.Range("F3").Value = 6397.836
$ws.Range("G3").Value = 5271.112
$ws.Range("H3").Value = 1126.725
$ws.Range("I3").Value = 1126.725
$ws.Range("K3").Value = 6620.928
$ws.Range("L3").Value = 5271.161
$ws.Range("M3").Value = 1349.767
$ws.Range("N3").Value = 1349.767
$ws.Range("B4").Value = 22.284
$ws.Range("C4").Value = 9.394
$ws.Range("D4").Value = 219.202
$ws.Range("F4").Value = 7281.11
$ws.Range("G4").Value = 5460.984
$ws.Range("H4").Value = 1820.126
$ws.Range("I4").Value = 1820.126
$ws.Range("K4").Value = 6616.0584
$ws.Range("L4").Value = 5460.984
$ws.Range("M4").Value = 1155.0744
$ws.Range("N4").Value = 1155.0744

$ws = $wb.Worksheets.Item("Fold_4")
$ws.Range("B2").Value = 22.51725
$ws.Range("C2").Value = 9.335249999999998
$ws.Range("D2").Value = 193.1642
$ws.Range("F2").Value = 6287.32935
$ws.Range("G2").Value = 5290.304150000001
$ws.Range("H2").Value = 997.0254000000001
$ws.Range("I2").Value = 997.0254000000001
$ws.Range("K2").Value = 6759.1448
$ws.Range("L2").Value = 5290.316
$ws.Range("M2").Value = 1468.8288
$ws.Range("N2").Value = 1468.8288
$ws.Range("B3").Value = 22.463
$ws.Range("C3").Value = 9.284000000000001
$ws.Range("D3").Value = 192.978
$ws.Range("F3").Value = 6362.588
$ws.Range("G3").Value = 5280.59
$ws.Range("H3").Value = 1081.998
$ws.Range("I3").Value = 1081.998
$ws.Range("K3").Value = 6760.028200000001
$ws.Range("L3").Value = 5280.656
$ws.Range("M3").Value = 1479.3722
$ws.Range("N3").Value = 1479.3722
$ws.Range("B4").Value = 26.357
$ws.Range("C4").Value = 10.383
$ws.Range("D4").Value = 209.797
$ws.Range("F4").Value = 7515.411
$ws.Range("G4").Value = 5900.331
$ws.Range("H4").Value = 1615.08
$ws.Range("I4").Value = 1615.08
$ws.Range("K4").Value = 6801.3754
$ws.Range("L4").Value = 5900.331
$ws.Range("M4").Value = 901.0444
$ws.Range("N4").Value = 901.0444

$ws = $wb.Worksheets.Item("Fold_5")
$ws.Range("B2").Value = 22.5829
$ws.Range("C2").Value = 9.088249999999999
$ws.Range("D2").Value = 196.0308
$ws.Range("F2").Value = 6380.1515
$ws.Range("G2").Value = 5309.8822
$ws.Range("H2").Value = 1070.2695
$ws.Range("I2").Value = 1070.2695
$ws.Range("K2").Value = 6377.329
$ws.Range("L2").Value = 5309.897
$ws.Range("M2").Value = 1067.432
$ws.Range("N2").Value = 1067.432
$ws.Range("B3").Value = 22.642
$ws.Range("C3").Value = 9.195
$ws.Range("D3").Value = 192.508
$ws.Range("F3").Value = 6458.032
$ws.Range("G3").Value = 5294.496
$ws.Range("H3").Value = 1163.536
$ws.Range("I3").Value = 1163.536
$ws.Range("K3").Value = 6376.9068
$ws.Range("L3").Value = 5294.425
$ws.Range("M3").Value = 1082.4818
$ws.Range("N3").Value = 1082.4818
$ws.Range("B4").Value = 26.357
$ws.Range("C4").Value = 10.383
$ws.Range("D4").Value = 209.797
$ws.Range("F4").Value = 7515.411
$ws.Range("G4").Value = 5900.331
$ws.Range("H4").Value = 1615.08
$ws.Range("I4").Value = 1615.08
$ws.Range("K4").Value = 6514.829400000001
$ws.Range("L4").Value = 5900.331
$ws.Range("M4").Value = 614.4984000000001
$ws.Range("N4").Value = 614.4984000000001

Write-Host "Updated all 5 Fold sheets with new plotting script values."